# Swap the species-observation data between row 7 and row 8, while
# leaving the shared/location columns (D, P, S, T, U, V, W, Y, AA, AD,
# AE, AG, AT, AW, AX, AY) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 7
$row2 = 8

# Columns whose values differ between the two rows and must be swapped.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"

    $val1 = $ws.Range($addr1).Value()
    $val2 = $ws.Range($addr2).Value()

    $ws.Range($addr1).Value = $val2
    $ws.Range($addr2).Value = $val1
}
